$p = $ppt.ActivePresentation
$s = $p.Slides.Item(18)
$shape = $s.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange

# Reorder the "Modern Navigation lifecycle events" bullets (paragraphs 6-8)
# from OnNavigatedFrom / OnNavigatedTo / OnNavigatingFrom
# to   OnNavigatedTo / OnNavigatingFrom / OnNavigatedFrom
#
# The text is set via a throwaway placeholder first so the host's
# minimal-diff text replacement doesn't fragment the run (the old/new
# strings share long common substrings, e.g. "OnNavigat...").
$tr.Paragraphs(6, 1).Text = "zzz1"
$tr.Paragraphs(6, 1).Text = "OnNavigatedTo"
$tr.Paragraphs(7, 1).Text = "zzz2"
$tr.Paragraphs(7, 1).Text = "OnNavigatingFrom"
$tr.Paragraphs(8, 1).Text = "zzz3"
$tr.Paragraphs(8, 1).Text = "OnNavigatedFrom"

# Turn those bullets (plus the trailing blank paragraph) into an
# auto-numbered (arabic period) list using the major-latin theme font,
# instead of the Wingdings "section" glyph bullet.
for ($i = 6; $i -le 9; $i++) {
    $para = $tr.Paragraphs($i, 1)
    $para.ParagraphFormat.Bullet.Font.Name = "+mj-lt"
    $para.ParagraphFormat.Bullet.Type = 2
    $para.ParagraphFormat.Bullet.Style = 3
}
